# Update the "netto" values (column B) for the first three cost positions
# on the "Gesamtinvestitionskosten" sheet from 12, 12, 7 to 111, 111, 111.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtinvestitionskosten")

$ws.Range("B2").Value = 111
$ws.Range("B3").Value = 111
$ws.Range("B4").Value = 111

$wb.Application.Calculate()
